$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 477, shifting existing rows 477-542 down to 478-543.
$ws.Rows.Item(477).Insert()

# Populate the newly inserted row 477 with the new data record.
$ws.Cells.Item(477, 1).Value = 3
$ws.Cells.Item(477, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(477, 3).Value = "Coquimbo"
$ws.Cells.Item(477, 4).Value = 45077
$ws.Cells.Item(477, 4).Style = $ws.Cells.Item(478, 4).Style
$ws.Cells.Item(477, 4).NumberFormat = $ws.Cells.Item(478, 4).NumberFormat
$ws.Cells.Item(477, 5).Value = 5
$ws.Cells.Item(477, 6).Value = 100112012
$ws.Cells.Item(477, 7).Value = "Espinaca"
$ws.Cells.Item(477, 8).Value = "Sin especificar"
$ws.Cells.Item(477, 9).Value = "Primera"
$ws.Cells.Item(477, 10).Value = 130
$ws.Cells.Item(477, 11).Value = 5000
$ws.Cells.Item(477, 12).Value = 5500
$ws.Cells.Item(477, 13).Value = 5231
$ws.Cells.Item(477, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(477, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(477, 16).Value = 1744
$ws.Cells.Item(477, 17).Value = 3
$ws.Cells.Item(477, 18).Value = "Hortaliza"
